# Unipept lookup sheet for 100 m sus fungi: add a new "Sheet3" listing the
# unique stripped peptides (incl. a newly-noticed AGQLEQIR) pulled from the
# noHi3 tab, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$peptides = @(
    "TATQKTVDGPSAK",
    "TTGIVLDSGDGVTH",
    "ELAEDGYSGVEVR",
    "DSYVGDEAQSKR",
    "QIVGDDLTVTNPK",
    "SGDSAIVK",
    "AGQLEQIR",
    "GDDLTVTNPK",
    "APANVTTEVK"
)

# Leave the last-used selections on the existing tabs as they were when the
# user tabbed away from them.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G11").Select() | Out-Null

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F10").Select() | Out-Null

# New sheet, appended after noHi3, carrying the stripped-peptide list.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

for ($i = 0; $i -lt $peptides.Length; $i++) {
    $ws3.Cells.Item($i + 1, 1).Value = $peptides[$i]
}

$ws3.Activate()
$ws3.Range("L18").Select() | Out-Null
